$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(132, 4).Value = "2021-11-18"
$ws.Cells.Item(132, 10).Value = 500
$ws.Cells.Item(132, 11).Value = 1000
$ws.Cells.Item(132, 12).Value = 1000
$ws.Cells.Item(132, 13).Value = 1000
$ws.Cells.Item(132, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(132, 15).Value = 'Región del Maule'
$ws.Cells.Item(132, 16).Value = 200
$ws.Cells.Item(132, 17).Value = 5

$ws.Cells.Item(133, 4).Value = "2021-09-07"
$ws.Cells.Item(133, 10).Value = 1200
$ws.Cells.Item(133, 11).Value = 1000
$ws.Cells.Item(133, 12).Value = 1000
$ws.Cells.Item(133, 13).Value = 1000
$ws.Cells.Item(133, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(133, 15).Value = 'Región del Maule'
$ws.Cells.Item(133, 16).Value = 200
$ws.Cells.Item(133, 17).Value = 5

$ws.Cells.Item(134, 4).Value = "2021-09-24"
$ws.Cells.Item(134, 10).Value = 1200
$ws.Cells.Item(134, 11).Value = 1200
$ws.Cells.Item(134, 12).Value = 1200
$ws.Cells.Item(134, 13).Value = 1200
$ws.Cells.Item(134, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(134, 15).Value = 'Región del Maule'
$ws.Cells.Item(134, 16).Value = 240
$ws.Cells.Item(134, 17).Value = 5

$ws.Cells.Item(135, 4).Value = "2021-02-18"
$ws.Cells.Item(135, 10).Value = 500
$ws.Cells.Item(135, 11).Value = 1000
$ws.Cells.Item(135, 12).Value = 1000
$ws.Cells.Item(135, 13).Value = 1000
$ws.Cells.Item(135, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(135, 15).Value = 'Región del Maule'
$ws.Cells.Item(135, 16).Value = 200
$ws.Cells.Item(135, 17).Value = 5

$ws.Cells.Item(136, 4).Value = "2021-05-07"
$ws.Cells.Item(136, 10).Value = 120
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = 9000
$ws.Cells.Item(136, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(136, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(136, 16).Value = 600
$ws.Cells.Item(136, 17).Value = 15

$ws.Cells.Item(137, 4).Value = "2021-05-07"
$ws.Cells.Item(137, 10).Value = 1000
$ws.Cells.Item(137, 11).Value = 1000
$ws.Cells.Item(137, 12).Value = 1000
$ws.Cells.Item(137, 13).Value = 1000
$ws.Cells.Item(137, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(137, 15).Value = 'Región del Maule'
$ws.Cells.Item(137, 16).Value = 200
$ws.Cells.Item(137, 17).Value = 5

$ws.Cells.Item(138, 4).Value = "2021-10-12"
$ws.Cells.Item(138, 10).Value = 1200
$ws.Cells.Item(138, 11).Value = 900
$ws.Cells.Item(138, 12).Value = 1000
$ws.Cells.Item(138, 13).Value = 950
$ws.Cells.Item(138, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(138, 15).Value = 'Región del Maule'
$ws.Cells.Item(138, 16).Value = 190
$ws.Cells.Item(138, 17).Value = 5

$ws.Cells.Item(139, 4).Value = "2021-02-02"
$ws.Cells.Item(139, 10).Value = 500
$ws.Cells.Item(139, 11).Value = 1000
$ws.Cells.Item(139, 12).Value = 1000
$ws.Cells.Item(139, 13).Value = 1000
$ws.Cells.Item(139, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(139, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(139, 16).Value = 200
$ws.Cells.Item(139, 17).Value = 5

$ws.Cells.Item(140, 4).Value = "2021-02-02"
$ws.Cells.Item(140, 10).Value = 500
$ws.Cells.Item(140, 11).Value = 1000
$ws.Cells.Item(140, 12).Value = 1000
$ws.Cells.Item(140, 13).Value = 1000
$ws.Cells.Item(140, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(140, 15).Value = 'Región del Maule'
$ws.Cells.Item(140, 16).Value = 200
$ws.Cells.Item(140, 17).Value = 5

$ws.Cells.Item(141, 4).Value = "2021-08-09"
$ws.Cells.Item(141, 10).Value = 500
$ws.Cells.Item(141, 11).Value = 1000
$ws.Cells.Item(141, 12).Value = 1000
$ws.Cells.Item(141, 13).Value = 1000
$ws.Cells.Item(141, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(141, 15).Value = 'Región del Maule'
$ws.Cells.Item(141, 16).Value = 200
$ws.Cells.Item(141, 17).Value = 5

$ws.Cells.Item(142, 4).Value = "2021-09-23"
$ws.Cells.Item(142, 10).Value = 500
$ws.Cells.Item(142, 11).Value = 1200
$ws.Cells.Item(142, 12).Value = 1200
$ws.Cells.Item(142, 13).Value = 1200
$ws.Cells.Item(142, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(142, 15).Value = 'Región del Maule'
$ws.Cells.Item(142, 16).Value = 240
$ws.Cells.Item(142, 17).Value = 5

$ws.Cells.Item(143, 4).Value = "2021-03-04"
$ws.Cells.Item(143, 10).Value = 500
$ws.Cells.Item(143, 11).Value = 1000
$ws.Cells.Item(143, 12).Value = 1000
$ws.Cells.Item(143, 13).Value = 1000
$ws.Cells.Item(143, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(143, 15).Value = 'Región del Maule'
$ws.Cells.Item(143, 16).Value = 200
$ws.Cells.Item(143, 17).Value = 5

$ws.Cells.Item(144, 4).Value = "2021-01-20"
$ws.Cells.Item(144, 10).Value = 250
$ws.Cells.Item(144, 11).Value = 900
$ws.Cells.Item(144, 12).Value = 900
$ws.Cells.Item(144, 13).Value = 900
$ws.Cells.Item(144, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(144, 15).Value = 'Región del Maule'
$ws.Cells.Item(144, 16).Value = 180
$ws.Cells.Item(144, 17).Value = 5

$ws.Cells.Item(145, 4).Value = "2021-07-29"
$ws.Cells.Item(145, 10).Value = 500
$ws.Cells.Item(145, 11).Value = 1000
$ws.Cells.Item(145, 12).Value = 1000
$ws.Cells.Item(145, 13).Value = 1000
$ws.Cells.Item(145, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(145, 15).Value = 'Región del Maule'
$ws.Cells.Item(145, 16).Value = 200
$ws.Cells.Item(145, 17).Value = 5

$ws.Cells.Item(146, 4).Value = "2021-04-09"
$ws.Cells.Item(146, 10).Value = 1200
$ws.Cells.Item(146, 11).Value = 1000
$ws.Cells.Item(146, 12).Value = 1000
$ws.Cells.Item(146, 13).Value = 1000
$ws.Cells.Item(146, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(146, 15).Value = 'Región del Maule'
$ws.Cells.Item(146, 16).Value = 200
$ws.Cells.Item(146, 17).Value = 5

$ws.Cells.Item(147, 4).Value = "2021-03-15"
$ws.Cells.Item(147, 10).Value = 500
$ws.Cells.Item(147, 11).Value = 800
$ws.Cells.Item(147, 12).Value = 1000
$ws.Cells.Item(147, 13).Value = 900
$ws.Cells.Item(147, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(147, 15).Value = 'Región del Maule'
$ws.Cells.Item(147, 16).Value = 180
$ws.Cells.Item(147, 17).Value = 5

$ws.Cells.Item(148, 4).Value = "2021-06-16"
$ws.Cells.Item(148, 10).Value = 100
$ws.Cells.Item(148, 11).Value = 1000
$ws.Cells.Item(148, 12).Value = 1000
$ws.Cells.Item(148, 13).Value = 1000
$ws.Cells.Item(148, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(148, 15).Value = 'Región del Maule'
$ws.Cells.Item(148, 16).Value = 200
$ws.Cells.Item(148, 17).Value = 5

$ws.Cells.Item(149, 4).Value = "2021-04-13"
$ws.Cells.Item(149, 10).Value = 1250
$ws.Cells.Item(149, 11).Value = 1000
$ws.Cells.Item(149, 12).Value = 1000
$ws.Cells.Item(149, 13).Value = 1000
$ws.Cells.Item(149, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(149, 15).Value = 'Región del Maule'
$ws.Cells.Item(149, 16).Value = 200
$ws.Cells.Item(149, 17).Value = 5

$ws.Cells.Item(150, 4).Value = "2021-03-02"
$ws.Cells.Item(150, 10).Value = 1200
$ws.Cells.Item(150, 11).Value = 1000
$ws.Cells.Item(150, 12).Value = 1000
$ws.Cells.Item(150, 13).Value = 1000
$ws.Cells.Item(150, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(150, 15).Value = 'Región del Maule'
$ws.Cells.Item(150, 16).Value = 200
$ws.Cells.Item(150, 17).Value = 5

$ws.Cells.Item(151, 4).Value = "2021-07-26"
$ws.Cells.Item(151, 10).Value = 500
$ws.Cells.Item(151, 11).Value = 1000
$ws.Cells.Item(151, 12).Value = 1000
$ws.Cells.Item(151, 13).Value = 1000
$ws.Cells.Item(151, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(151, 15).Value = 'Región del Maule'
$ws.Cells.Item(151, 16).Value = 200
$ws.Cells.Item(151, 17).Value = 5

$ws.Cells.Item(152, 4).Value = "2020-12-04"
$ws.Cells.Item(152, 10).Value = 120
$ws.Cells.Item(152, 11).Value = 8000
$ws.Cells.Item(152, 12).Value = 8000
$ws.Cells.Item(152, 13).Value = 8000
$ws.Cells.Item(152, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(152, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(152, 16).Value = 533
$ws.Cells.Item(152, 17).Value = 15

$ws.Cells.Item(153, 4).Value = "2020-12-04"
$ws.Cells.Item(153, 10).Value = 1200
$ws.Cells.Item(153, 11).Value = 850
$ws.Cells.Item(153, 12).Value = 1000
$ws.Cells.Item(153, 13).Value = 925
$ws.Cells.Item(153, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(153, 15).Value = 'Región del Maule'
$ws.Cells.Item(153, 16).Value = 185
$ws.Cells.Item(153, 17).Value = 5

$ws.Cells.Item(154, 4).Value = "2021-06-29"
$ws.Cells.Item(154, 10).Value = 110
$ws.Cells.Item(154, 11).Value = 9000
$ws.Cells.Item(154, 12).Value = 9000
$ws.Cells.Item(154, 13).Value = 9000
$ws.Cells.Item(154, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(154, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(154, 16).Value = 600
$ws.Cells.Item(154, 17).Value = 15

$ws.Cells.Item(155, 4).Value = "2021-06-29"
$ws.Cells.Item(155, 10).Value = 1200
$ws.Cells.Item(155, 11).Value = 1000
$ws.Cells.Item(155, 12).Value = 1000
$ws.Cells.Item(155, 13).Value = 1000
$ws.Cells.Item(155, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(155, 15).Value = 'Región del Maule'
$ws.Cells.Item(155, 16).Value = 200
$ws.Cells.Item(155, 17).Value = 5

$ws.Cells.Item(156, 4).Value = "2021-10-05"
$ws.Cells.Item(156, 10).Value = 1400
$ws.Cells.Item(156, 11).Value = 1000
$ws.Cells.Item(156, 12).Value = 1000
$ws.Cells.Item(156, 13).Value = 1000
$ws.Cells.Item(156, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(156, 15).Value = 'Región del Maule'
$ws.Cells.Item(156, 16).Value = 200
$ws.Cells.Item(156, 17).Value = 5

$ws.Cells.Item(157, 4).Value = "2020-12-07"
$ws.Cells.Item(157, 10).Value = 500
$ws.Cells.Item(157, 11).Value = 800
$ws.Cells.Item(157, 12).Value = 800
$ws.Cells.Item(157, 13).Value = 800
$ws.Cells.Item(157, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(157, 15).Value = 'Región del Maule'
$ws.Cells.Item(157, 16).Value = 160
$ws.Cells.Item(157, 17).Value = 5

$ws.Cells.Item(158, 4).Value = "2021-08-13"
$ws.Cells.Item(158, 10).Value = 1200
$ws.Cells.Item(158, 11).Value = 1000
$ws.Cells.Item(158, 12).Value = 1000
$ws.Cells.Item(158, 13).Value = 1000
$ws.Cells.Item(158, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(158, 15).Value = 'Región del Maule'
$ws.Cells.Item(158, 16).Value = 200
$ws.Cells.Item(158, 17).Value = 5

$ws.Cells.Item(159, 4).Value = "2021-02-15"
$ws.Cells.Item(159, 10).Value = 500
$ws.Cells.Item(159, 11).Value = 1000
$ws.Cells.Item(159, 12).Value = 1000
$ws.Cells.Item(159, 13).Value = 1000
$ws.Cells.Item(159, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(159, 15).Value = 'Región del Maule'
$ws.Cells.Item(159, 16).Value = 200
$ws.Cells.Item(159, 17).Value = 5

$ws.Cells.Item(160, 4).Value = "2021-02-12"
$ws.Cells.Item(160, 10).Value = 1200
$ws.Cells.Item(160, 11).Value = 1000
$ws.Cells.Item(160, 12).Value = 1000
$ws.Cells.Item(160, 13).Value = 1000
$ws.Cells.Item(160, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(160, 15).Value = 'Región del Maule'
$ws.Cells.Item(160, 16).Value = 200
$ws.Cells.Item(160, 17).Value = 5

$ws.Cells.Item(161, 4).Value = "2021-01-26"
$ws.Cells.Item(161, 10).Value = 1000
$ws.Cells.Item(161, 11).Value = 900
$ws.Cells.Item(161, 12).Value = 900
$ws.Cells.Item(161, 13).Value = 900
$ws.Cells.Item(161, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(161, 15).Value = 'Región del Maule'
$ws.Cells.Item(161, 16).Value = 180
$ws.Cells.Item(161, 17).Value = 5

$ws.Cells.Item(162, 4).Value = "2021-09-09"
$ws.Cells.Item(162, 10).Value = 500
$ws.Cells.Item(162, 11).Value = 1000
$ws.Cells.Item(162, 12).Value = 1000
$ws.Cells.Item(162, 13).Value = 1000
$ws.Cells.Item(162, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(162, 15).Value = 'Región del Maule'
$ws.Cells.Item(162, 16).Value = 200
$ws.Cells.Item(162, 17).Value = 5

$ws.Cells.Item(163, 4).Value = "2021-06-15"
$ws.Cells.Item(163, 10).Value = 110
$ws.Cells.Item(163, 11).Value = 10000
$ws.Cells.Item(163, 12).Value = 10000
$ws.Cells.Item(163, 13).Value = 10000
$ws.Cells.Item(163, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(163, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(163, 16).Value = 667
$ws.Cells.Item(163, 17).Value = 15

$ws.Cells.Item(164, 4).Value = "2021-06-15"
$ws.Cells.Item(164, 10).Value = 1200
$ws.Cells.Item(164, 11).Value = 1000
$ws.Cells.Item(164, 12).Value = 1000
$ws.Cells.Item(164, 13).Value = 1000
$ws.Cells.Item(164, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(164, 15).Value = 'Región del Maule'
$ws.Cells.Item(164, 16).Value = 200
$ws.Cells.Item(164, 17).Value = 5

$ws.Cells.Item(165, 4).Value = "2021-01-14"
$ws.Cells.Item(165, 10).Value = 750
$ws.Cells.Item(165, 11).Value = 900
$ws.Cells.Item(165, 12).Value = 1000
$ws.Cells.Item(165, 13).Value = 967
$ws.Cells.Item(165, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(165, 15).Value = 'Región del Maule'
$ws.Cells.Item(165, 16).Value = 193
$ws.Cells.Item(165, 17).Value = 5

$ws.Cells.Item(166, 4).Value = "2020-12-11"
$ws.Cells.Item(166, 10).Value = 1200
$ws.Cells.Item(166, 11).Value = 850
$ws.Cells.Item(166, 12).Value = 850
$ws.Cells.Item(166, 13).Value = 850
$ws.Cells.Item(166, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(166, 15).Value = 'Región del Maule'
$ws.Cells.Item(166, 16).Value = 170
$ws.Cells.Item(166, 17).Value = 5

$ws.Cells.Item(167, 4).Value = "2021-04-15"
$ws.Cells.Item(167, 10).Value = 500
$ws.Cells.Item(167, 11).Value = 1000
$ws.Cells.Item(167, 12).Value = 1000
$ws.Cells.Item(167, 13).Value = 1000
$ws.Cells.Item(167, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(167, 15).Value = 'Región del Maule'
$ws.Cells.Item(167, 16).Value = 200
$ws.Cells.Item(167, 17).Value = 5
